$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("anon_post")

# Insert a new row before row 8, shifting existing rows (old row 8 "Q50" etc.) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new "Q54 / gender" lookup entry.
# Column order of entry (A, then C, then B) reproduces the shared-string insertion
# order required (Q54, "What is your gender?", "gender: 1=woman, 2=man, 3=other").
$ws.Range("A8").Value2 = "Q54"
$ws.Range("C8").Value2 = "What is your gender?"
$ws.Range("B8").Value2 = "gender: 1=woman, 2=man, 3=other"

# Update the hidden autofilter defined name so it still spans the data block now that
# it has grown by one row.
$filterName = $wb.Names.Item("anon_post!_FilterDatabase")
$filterName.RefersTo = "=anon_post!`$B`$27:`$C`$111"

# Match the saved selection/active cell on the anon_post sheet.
$ws.Activate()
$ws.Range("B8").Select()
